$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the obsolete year rows (2000年, 2002年, 2005年, 2007年) ---
# These are currently rows 2-5; deleting them shifts the existing
# 2010年/2012年/2015年/2017年 rows (6-9) up to become rows 2-5.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# --- Append the new 2020年 row as row 6 ---
# Copy formatting from the cell above (A5, 2017年) so the new year
# cell matches the existing style (bold/centered with border).
$ws.Cells.Item(5,1).Copy()
$ws.Cells.Item(6,1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(6,1).Value = "2020年"
$ws.Cells.Item(6,2).Value = 1945824.44731674
$ws.Cells.Item(6,3).Value = 5234595.238605
$ws.Cells.Item(6,4).Value = 17577339.566559
$ws.Cells.Item(6,5).Value = ""
$ws.Cells.Item(6,6).Value = 4243704003.37183
$ws.Cells.Item(6,7).Value = 2707176325.74967
$ws.Cells.Item(6,8).Value = ""
$ws.Cells.Item(6,9).Value = ""
$ws.Cells.Item(6,10).Value = 783222749.856351
$ws.Cells.Item(6,11).Value = ""
$ws.Cells.Item(6,12).Value = ""
$ws.Cells.Item(6,13).Value = ""
$ws.Cells.Item(6,14).Value = ""
$ws.Cells.Item(6,15).Value = ""
$ws.Cells.Item(6,16).Value = 25162931.8981722
$ws.Cells.Item(6,17).Value = ""
$ws.Cells.Item(6,18).Value = ""
$ws.Cells.Item(6,19).Value = ""

# --- Clear anything below the new last data row (row 6) ---
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -gt 6) {
    $ws.Range($ws.Cells.Item(7,1), $ws.Cells.Item($lastRow,19)).Delete(-4162)
}
